$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '26.220.14'
Set-TextCell 'E2' '  -0.67%  '
Set-TextCell 'D3' '1.681.25'
Set-TextCell 'E3' '  -0.65%  '
Set-TextCell 'E4' '  -0.64%  '
Set-TextCell 'D5' '211.28'
Set-TextCell 'E5' '  -3.52%  '
Set-TextCell 'D6' '0.5319'
Set-TextCell 'E6' '  -4.14%  '
Set-TextCell 'D7' '1.004'
Set-TextCell 'E7' '  -0.72%  '
Set-TextCell 'D8' '0.2689'
Set-TextCell 'E8' '  -1.17%  '
Set-TextCell 'D9' '0.06312'
Set-TextCell 'E9' '  -2.70%  '
Set-TextCell 'D10' '21.41'
Set-TextCell 'E10' '  -3.40%  '
Set-TextCell 'D11' '0.07529'
Set-TextCell 'E11' '  -1.85%  '
Set-TextCell 'D12' '1.687.26'
Set-TextCell 'E12' '  -0.21%  '
Set-TextCell 'D13' '4.498'
Set-TextCell 'E13' '  -1.41%  '
Set-TextCell 'D14' '0.5682'
Set-TextCell 'E14' '  -2.52%  '
Set-TextCell 'D15' '0.000008169'
Set-TextCell 'E15' '  -3.50%  '
Set-TextCell 'D16' '66.64'
Set-TextCell 'E16' '  +2.00%  '
Set-TextCell 'D17' '26.243.10'
Set-TextCell 'E17' '  -0.98%  '
Set-TextCell 'E18' '  -0.56%  '
Set-TextCell 'D19' '4.859'
Set-TextCell 'E19' '  -2.20%  '
Set-TextCell 'E20' '  -4.08%  '
Set-TextCell 'D21' '189.21'
Set-TextCell 'E21' '  -0.63%  '
Set-TextCell 'D22' '6.226'
Set-TextCell 'E22' '  -0.44%  '
Set-TextCell 'D23' '1.005'
Set-TextCell 'E23' '  -0.62%  '
Set-TextCell 'D24' '147.87'
Set-TextCell 'E24' '  -1.41%  '
Set-TextCell 'D25' '0.1270'
Set-TextCell 'E25' '  -2.91%  '
Set-TextCell 'D26' '7.640'
Set-TextCell 'E26' '  -3.42%  '
Set-TextCell 'D27' '15.94'
Set-TextCell 'E27' '  +1.15%  '
Set-TextCell 'D28' '0.06465'
Set-TextCell 'E28' '  +2.11%  '
Set-TextCell 'D29' '1.342'
Set-TextCell 'E29' '  -5.65%  '
Set-TextCell 'D30' '1.285'
Set-TextCell 'E30' '  -3.48%  '
Set-TextCell 'D31' '3.534'
Set-TextCell 'E31' '  -1.76%  '
Set-TextCell 'D32' '3.489'
Set-TextCell 'E32' '  -3.02%  '
Set-TextCell 'D33' '1.661'
Set-TextCell 'E33' '  -1.02%  '
Set-TextCell 'D34' '1.013'
Set-TextCell 'E34' '  -3.01%  '
Set-TextCell 'D35' '0.6114'
Set-TextCell 'E35' '  -1.66%  '
Set-TextCell 'D36' '2.413'
Set-TextCell 'E36' '  +0.29%  '
Set-TextCell 'D37' '2.718'
Set-TextCell 'E37' '  -0.22%  '
Set-TextCell 'D38' '6.193'
Set-TextCell 'E38' '  -0.68%  '
Set-TextCell 'D39' '0.01622'
Set-TextCell 'E39' '  -1.26%  '
Set-TextCell 'D40' '1.105.18'
Set-TextCell 'E40' '  -1.79%  '
Set-TextCell 'D41' '0.8679'
Set-TextCell 'E41' '  -1.56%  '
Set-TextCell 'E42' '  -0.99%  '
Set-TextCell 'D43' '100.13'
Set-TextCell 'E43' '  -0.74%  '
Set-TextCell 'D44' '1.831.99'
Set-TextCell 'E44' '  -0.65%  '
Set-TextCell 'D45' '0.00000000107'
Set-TextCell 'E45' '  -2.76%  '
Set-TextCell 'D46' '56.95'
Set-TextCell 'E46' '  -1.09%  '
Set-TextCell 'E47' '  -0.27%  '
Set-TextCell 'B48' 'EnergySwap'
Set-TextCell 'C48' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D48' '8.038'
Set-TextCell 'E48' '  -2.47%  '
Set-TextCell 'B49' 'Cronos'
Set-TextCell 'C49' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 'D49' '0.05271'
Set-TextCell 'E49' '  -0.14%  '
Set-TextCell 'D50' '0.4276'
Set-TextCell 'E50' '  -0.61%  '
Set-TextCell 'D51' '5.974'
Set-TextCell 'E51' '  -1.75%  '
